{"js": "// Update the division problems in the table to the new set of values.\n// Each old expression is unique within the document, so an exact,\n// case-sensitive search-and-replace per pair is safe and unambiguous.\nconst replacements = [\n  [\"829\u00f76=\", \"532\u00f79=\"],\n  [\"603\u00f76=\", \"909\u00f77=\"],\n  [\"336\u00f73=\", \"471\u00f77=\"],\n  [\"591\u00f72=\", \"880\u00f78=\"],\n  [\"849\u00f74=\", \"835\u00f77=\"],\n  [\"141\u00f79=\", \"327\u00f72=\"],\n  [\"651\u00f75=\", \"649\u00f72=\"],\n  [\"506\u00f78=\", \"948\u00f79=\"],\n  [\"508\u00f76=\", \"214\u00f72=\"],\n  [\"677\u00f74=\", \"768\u00f77=\"],\n  [\"489\u00f75=\", \"279\u00f73=\"],\n  [\"213\u00f78=\", \"712\u00f78=\"],\n  [\"707\u00f72=\", \"476\u00f77=\"],\n  [\"679\u00f76=\", \"646\u00f79=\"],\n  [\"601\u00f78=\", \"530\u00f73=\"],\n  [\"819\u00f75=\", \"214\u00f79=\"],\n  [\"146\u00f79=\", \"148\u00f79=\"],\n  [\"533\u00f72=\", \"750\u00f76=\"],\n  [\"851\u00f79=\", \"248\u00f73=\"],\n  [\"547\u00f72=\", \"978\u00f74=\"],\n  [\"377\u00f74=\", \"914\u00f75=\"],\n  [\"978\u00f73=\", \"884\u00f72=\"],\n  [\"814\u00f79=\", \"193\u00f72=\"],\n  [\"402\u00f75=\", \"276\u00f73=\"],\n  [\"962\u00f72=\", \"288\u00f75=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the division problems in the table to the new set of values.\n# Each old expression appears exactly once in the document, so a plain\n# Find/Replace (wdReplaceAll, exact text) per pair is unambiguous.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"829\u00f76=\", \"532\u00f79=\"),\n    @(\"603\u00f76=\", \"909\u00f77=\"),\n    @(\"336\u00f73=\", \"471\u00f77=\"),\n    @(\"591\u00f72=\", \"880\u00f78=\"),\n    @(\"849\u00f74=\", \"835\u00f77=\"),\n    @(\"141\u00f79=\", \"327\u00f72=\"),\n    @(\"651\u00f75=\", \"649\u00f72=\"),\n    @(\"506\u00f78=\", \"948\u00f79=\"),\n    @(\"508\u00f76=\", \"214\u00f72=\"),\n    @(\"677\u00f74=\", \"768\u00f77=\"),\n    @(\"489\u00f75=\", \"279\u00f73=\"),\n    @(\"213\u00f78=\", \"712\u00f78=\"),\n    @(\"707\u00f72=\", \"476\u00f77=\"),\n    @(\"679\u00f76=\", \"646\u00f79=\"),\n    @(\"601\u00f78=\", \"530\u00f73=\"),\n    @(\"819\u00f75=\", \"214\u00f79=\"),\n    @(\"146\u00f79=\", \"148\u00f79=\"),\n    @(\"533\u00f72=\", \"750\u00f76=\"),\n    @(\"851\u00f79=\", \"248\u00f73=\"),\n    @(\"547\u00f72=\", \"978\u00f74=\"),\n    @(\"377\u00f74=\", \"914\u00f75=\"),\n    @(\"978\u00f73=\", \"884\u00f72=\"),\n    @(\"814\u00f79=\", \"193\u00f72=\"),\n    @(\"402\u00f75=\", \"276\u00f73=\"),\n    @(\"962\u00f72=\", \"288\u00f75=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1  # wdFindContinue\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}\n"}
